$d = $word.ActiveDocument

# --- 1. "There is no 'Please select' line in select controls in Edit mode
#        (New mode is fine)" -> split the sentence: keep the first clause as
#        is and extend the trailing clause with the extra detail about the
#        manager field.
$old   = "There is no ‘Please select’ line in select controls in Edit mode (New mode is fine)"
$part1 = "There is no ‘Please select’ line in select controls"
$part2 = " in Edit mode (New mode is fine, manager field in user is fine, too)"

[void]$d.Content.Find.Execute($old, $true, $false, $false, $false, $false, $true, 1, $false, $part1, 2)

$targetPara = $null
for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
    $paraText = $d.Paragraphs($i).Range.Text.TrimEnd([char]13)
    if ($paraText -eq $part1) {
        $targetPara = $d.Paragraphs($i)
        break
    }
}

if ($targetPara -eq $null) {
    throw "Could not locate the 'Please select' bullet paragraph after the Find/Replace"
}

$tailRange = $targetPara.Range
$tailRange.Collapse(0)
[void]$tailRange.MoveEnd(1, -1)
$tailRange.InsertAfter($part2)

# --- 2. Add a new bullet right after it for the new ToDo item about guarding
#        selection variables against nil (e.g. the new accountant field on
#        company).
[void]$targetPara.Range.InsertParagraphAfter()
$newPara = $targetPara.Next()
$newPara.Range.Text = "Check selection variables for not nil (like accountant in company)"
